# Insert a new data row at row 62 (pushing existing rows 62:180 down to 63:181)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 62:180 down by one row -> rows 63:181
$ws.Rows("62:62").Insert()

# Populate the newly-inserted row 62 with the new record.
$ws.Range("A62").Value = 4
$ws.Range("B62").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C62").Value = "Los Lagos"
$ws.Range("D62").Value = 44540
$ws.Range("E62").Value = 10
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100103
$ws.Range("H62").Value = "Frutos de hueso (carozo)"
$ws.Range("I62").Value = 100103004
$ws.Range("J62").Value = "Durazno"
$ws.Range("K62").Value = "Florida King"
$ws.Range("L62").Value = "Tercera"
$ws.Range("M62").Value = 800
$ws.Range("N62").Value = 16000
$ws.Range("O62").Value = 16000
$ws.Range("P62").Value = 16000
$ws.Range("Q62").Value = "$/caja 15 kilos empedrada"
$ws.Range("R62").Value = "Provincia de Limarí"
$ws.Range("S62").Value = 1067
$ws.Range("T62").Value = 15
